# Auto-generated edit script applying market-price / profit-column updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (Pandaemonium_Profits data).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 38500
$ws.Range("J3").Value = 38500
$ws.Range("L3").Value = 38500
$ws.Range("N3").Value = -38728
$ws.Range("H19").Value = 847.6
$ws.Range("I19").Value = 530
$ws.Range("K19").Value = 530
$ws.Range("M19").Value = -355
$ws.Range("H32").Value = 725.1667
$ws.Range("I32").Value = 717
$ws.Range("J32").Value = 733.3333
$ws.Range("K32").Value = 717
$ws.Range("L32").Value = 733.3333
$ws.Range("M32").Value = -391
$ws.Range("N32").Value = -1385.3333
$ws.Range("H43").Value = 785.4286
$ws.Range("I43").Value = 750
$ws.Range("J43").Value = 998
$ws.Range("K43").Value = 750
$ws.Range("L43").Value = 998
$ws.Range("M43").Value = -681
$ws.Range("N43").Value = -1136
$ws.Range("H51").Value = 2571.1428
$ws.Range("I51").Value = 3000
$ws.Range("J51").Value = 2499.6667
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 2499.6667
$ws.Range("M51").Value = -2516
$ws.Range("N51").Value = -3467.6667
$ws.Range("H62").Value = 4127.0454
$ws.Range("I62").Value = 2444.75
$ws.Range("J62").Value = 20950
$ws.Range("K62").Value = 2444.75
$ws.Range("L62").Value = 20950
$ws.Range("M62").Value = -1820.75
$ws.Range("N62").Value = -22198
$ws.Range("H65").Value = 4127.0454
$ws.Range("I65").Value = 2444.75
$ws.Range("J65").Value = 20950
$ws.Range("K65").Value = 12223.75
$ws.Range("L65").Value = 104750
$ws.Range("M65").Value = -9103.75
$ws.Range("N65").Value = -110990
$ws.Range("H69").Value = 4117.5557
$ws.Range("I69").Value = 6000
$ws.Range("J69").Value = 3579.7144
$ws.Range("K69").Value = 18000
$ws.Range("L69").Value = 10739.1432
$ws.Range("M69").Value = -17126
$ws.Range("N69").Value = -12487.1432
$ws.Range("H72").Value = 4117.5557
$ws.Range("I72").Value = 6000
$ws.Range("J72").Value = 3579.7144
$ws.Range("K72").Value = 54000
$ws.Range("L72").Value = 32217.4296
$ws.Range("M72").Value = -49632
$ws.Range("N72").Value = -40953.4296
$ws.Range("H74").Value = 4988.3076
$ws.Range("I74").Value = 4283.1665
$ws.Range("J74").Value = 5592.7144
$ws.Range("K74").Value = 4283.1665
$ws.Range("L74").Value = 5592.7144
$ws.Range("M74").Value = -3347.1665
$ws.Range("N74").Value = -7464.7144
$ws.Range("H77").Value = 4988.3076
$ws.Range("I77").Value = 4283.1665
$ws.Range("J77").Value = 5592.7144
$ws.Range("K77").Value = 21415.8325
$ws.Range("L77").Value = 27963.572
$ws.Range("M77").Value = -16735.8325
$ws.Range("N77").Value = -37323.572
$ws.Range("H102").Value = 38500
$ws.Range("J102").Value = 38500
$ws.Range("L102").Value = 38500
$ws.Range("N102").Value = -44990
$ws.Range("H116").Value = 2493.25
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 2789.2
$ws.Range("K116").Value = 2000
$ws.Range("L116").Value = 2789.2
$ws.Range("M116").Value = 1442
$ws.Range("N116").Value = -9673.200000000001
$ws.Range("H138").Value = 11056.513
$ws.Range("I138").Value = 3499.5
$ws.Range("J138").Value = 13006.71
$ws.Range("K138").Value = 10498.5
$ws.Range("L138").Value = 39020.13
$ws.Range("M138").Value = -5358.5
$ws.Range("N138").Value = -49300.13

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 407.2857
$ws.Range("I5").Value = 330.2
$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 330.2
$ws.Range("L5").Value = 600
$ws.Range("M5").Value = -218.2
$ws.Range("N5").Value = -824
$ws.Range("H63").Value = 3186.25
$ws.Range("I63").Value = 2377.5
$ws.Range("J63").Value = 3995
$ws.Range("K63").Value = 2377.5
$ws.Range("L63").Value = 3995
$ws.Range("M63").Value = -1691.5
$ws.Range("N63").Value = -5367
$ws.Range("H66").Value = 3186.25
$ws.Range("I66").Value = 2377.5
$ws.Range("J66").Value = 3995
$ws.Range("K66").Value = 11887.5
$ws.Range("L66").Value = 19975
$ws.Range("M66").Value = -8455.5
$ws.Range("N66").Value = -26839
$ws.Range("H132").Value = 6364.304
$ws.Range("I132").Value = 7741.6875
$ws.Range("J132").Value = 3216
$ws.Range("K132").Value = 23225.0625
$ws.Range("L132").Value = 9648
$ws.Range("M132").Value = -20695.0625
$ws.Range("N132").Value = -14708

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 407.2857
$ws.Range("I4").Value = 330.2
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 330.2
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = -215.2
$ws.Range("N4").Value = -830
$ws.Range("H22").Value = 293.75
$ws.Range("I22").Value = 299.7143
$ws.Range("J22").Value = 252
$ws.Range("K22").Value = 299.7143
$ws.Range("L22").Value = 252
$ws.Range("M22").Value = -126.7143
$ws.Range("N22").Value = -598
$ws.Range("H86").Value = 1803.2808
$ws.Range("I86").Value = 1796.1538
$ws.Range("K86").Value = 1796.1538
$ws.Range("M86").Value = -673.1538
$ws.Range("H89").Value = 1803.2808
$ws.Range("I89").Value = 1796.1538
$ws.Range("K89").Value = 8980.769
$ws.Range("M89").Value = -3364.769
$ws.Range("H94").Value = 908.4761999999999
$ws.Range("I94").Value = 1048.625
$ws.Range("J94").Value = 460
$ws.Range("K94").Value = 1048.625
$ws.Range("L94").Value = 460
$ws.Range("M94").Value = -597.625
$ws.Range("N94").Value = -1362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 350000670
$ws.Range("I6").Value = 350000670
$ws.Range("K6").Value = 350000670
$ws.Range("M6").Value = -350000557
$ws.Range("H134").Value = 3006.4666
$ws.Range("I134").Value = 2393.1
$ws.Range("K134").Value = 7179.299999999999
$ws.Range("M134").Value = -4644.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 831.8889
$ws.Range("I114").Value = 552
$ws.Range("J114").Value = 971.8333
$ws.Range("K114").Value = 1656
$ws.Range("L114").Value = 2915.4999
$ws.Range("M114").Value = 1598
$ws.Range("N114").Value = -9423.499899999999
$ws.Range("H117").Value = 5033.3335
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 5033.3335
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 15100.0005
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -21984.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 21333.334
$ws.Range("I80").Value = 50000
$ws.Range("J80").Value = 7000
$ws.Range("K80").Value = 50000
$ws.Range("L80").Value = 7000
$ws.Range("M80").Value = -49002
$ws.Range("N80").Value = -8996
$ws.Range("H83").Value = 21333.334
$ws.Range("I83").Value = 50000
$ws.Range("J83").Value = 7000
$ws.Range("K83").Value = 250000
$ws.Range("L83").Value = 35000
$ws.Range("M83").Value = -245008
$ws.Range("N83").Value = -44984
$ws.Range("H126").Value = 2680.9546
$ws.Range("I126").Value = 1806.75
$ws.Range("J126").Value = 3730
$ws.Range("K126").Value = 5420.25
$ws.Range("L126").Value = 11190
$ws.Range("M126").Value = -2950.25
$ws.Range("N126").Value = -16130
$ws.Range("H132").Value = 8379.706
$ws.Range("I132").Value = 11292.8
$ws.Range("J132").Value = 4218.143
$ws.Range("K132").Value = 33878.39999999999
$ws.Range("L132").Value = 12654.429
$ws.Range("M132").Value = -31348.39999999999
$ws.Range("N132").Value = -17714.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2000.5
$ws.Range("I22").Value = 3001
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 3001
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -2706
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 2000.5
$ws.Range("I27").Value = 3001
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 3001
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -2894
$ws.Range("N27").Value = -1214
$ws.Range("H40").Value = 3520.5
$ws.Range("I40").Value = 2537.875
$ws.Range("K40").Value = 2537.875
$ws.Range("M40").Value = -2401.875
$ws.Range("H46").Value = 1066.3334
$ws.Range("J46").Value = 679.6
$ws.Range("L46").Value = 679.6
$ws.Range("N46").Value = -1055.6
$ws.Range("H100").Value = 3437.5
$ws.Range("I100").Value = 3583.3333
$ws.Range("K100").Value = 3583.3333
$ws.Range("M100").Value = -3042.3333
$ws.Range("H140").Value = 53543.145
$ws.Range("J140").Value = 53543.145
$ws.Range("L140").Value = 53543.145
$ws.Range("N140").Value = -63903.145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H122").Value = 8814.5
$ws.Range("I122").Value = 1728.2858
$ws.Range("K122").Value = 5184.857400000001
$ws.Range("M122").Value = -2734.857400000001
$ws.Range("H136").Value = 3766.2896
$ws.Range("I136").Value = 2775.35
$ws.Range("K136").Value = 8326.049999999999
$ws.Range("M136").Value = -5776.049999999999

